# "Switch all DK to ENG"
# Translate Danish labels to English across the input sheets, change the
# TRUE/FALSE "scenario active" flags on "Field app scenarios" into explicit
# =TRUE()/=FALSE() formulas, and leave the workbook with the "Other inputs"
# sheet active/selected (matching the saved state captured in the diff).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Raw slurry comp" - translate manure source names
# ---------------------------------------------------------------------
$wsRaw = $wb.Worksheets.Item("Raw slurry comp")
$wsRaw.Range("A2").Value = "Pig"
$wsRaw.Range("A3").Value = "Cattle"
$wsRaw.Range("A4").Value = "Digestate"

# ---------------------------------------------------------------------
# 2) "Application climate" - translate month / season names
# ---------------------------------------------------------------------
$wsClimate = $wb.Worksheets.Item("Application climate")
$wsClimate.Range("A2").Value = "March"
$wsClimate.Range("A4").Value = "May"
$wsClimate.Range("A5").Value = "Summer"
$wsClimate.Range("A6").Value = "Autumn"

# ---------------------------------------------------------------------
# 3) "Storage EFs" - same manure source translations as "Raw slurry comp"
# ---------------------------------------------------------------------
$wsStorage = $wb.Worksheets.Item("Storage EFs")
$wsStorage.Range("A2").Value = "Pig"
$wsStorage.Range("A3").Value = "Cattle"
$wsStorage.Range("A4").Value = "Digestate"

# ---------------------------------------------------------------------
# 4) "Field app scenarios" - turn the literal TRUE/FALSE values in column B
#    into explicit =TRUE()/=FALSE() formulas (same displayed result).
# ---------------------------------------------------------------------
$wsScenarios = $wb.Worksheets.Item("Field app scenarios")
foreach ($r in 2..3) {
    $wsScenarios.Range("B$r").Formula = "=FALSE()"
}
foreach ($r in 4..19) {
    $wsScenarios.Range("B$r").Formula = "=TRUE()"
}

# ---------------------------------------------------------------------
# 5) "Other inputs" - nudge the style on the header/value cells (re-apply
#    the Normal style) so the sheet carries its own explicit cell format.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("Other inputs")
$wsOther.Range("A1:B2").Style = "Normal"

# ---------------------------------------------------------------------
# 6) Restore / update the per-sheet selections left behind by the edit,
#    finishing on "Other inputs" so it becomes the active (saved) sheet.
# ---------------------------------------------------------------------
$wsScenarios.Range("A2").Select()
$wsRaw.Range("A5").Select()
$wsClimate.Range("A7").Select()
$wsStorage.Range("A5").Select()
$wsOther.Range("C8").Select()
